$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "29.221.47"
$ws.Cells.Item(2, 5).Value = "  -0.66%  "
$ws.Cells.Item(3, 4).Value = "1.861.58"
$ws.Cells.Item(3, 5).Value = "  -0.83%  "
$ws.Cells.Item(4, 4).NumberFormat = "@"
$ws.Cells.Item(4, 4).Value = "0.9992"
$ws.Cells.Item(4, 5).Value = "  -0.06%  "
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "0.7140"
$ws.Cells.Item(5, 5).Value = "  -0.79%  "
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = "240.57"
$ws.Cells.Item(6, 5).Value = "  +0.21%  "
$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = "1.0000"
$ws.Cells.Item(7, 5).Value = "  +0.00%  "
$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = "0.3089"
$ws.Cells.Item(8, 5).Value = "  -0.36%  "
$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = "0.07710"
$ws.Cells.Item(9, 5).Value = "  -1.42%  "
$ws.Cells.Item(10, 5).Value = "  +0.17%  "
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = "0.08309"
$ws.Cells.Item(11, 5).Value = "  +0.61%  "
$ws.Cells.Item(12, 4).Value = "1.880.30"
$ws.Cells.Item(12, 5).Value = "  +0.29%  "
$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = "0.7174"
$ws.Cells.Item(13, 5).Value = "  -1.34%  "
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = "5.214"
$ws.Cells.Item(14, 5).Value = "  -1.29%  "
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = "90.88"
$ws.Cells.Item(15, 5).Value = "  -0.40%  "
$ws.Cells.Item(16, 4).Value = "29.245.37"
$ws.Cells.Item(16, 5).Value = "  -0.34%  "
$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = "6.004"
$ws.Cells.Item(17, 5).Value = "  +1.53%  "
$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = "243.84"
$ws.Cells.Item(18, 5).Value = "  -0.47%  "
$ws.Cells.Item(19, 4).Value = "2.151.37"
$ws.Cells.Item(19, 5).Value = "  +2.56%  "
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = "0.000007804"
$ws.Cells.Item(20, 5).Value = "  -1.25%  "
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = "13.17"
$ws.Cells.Item(21, 5).Value = "  -1.03%  "
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = "0.9995"
$ws.Cells.Item(22, 5).Value = "  +0.01%  "
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = "7.958"
$ws.Cells.Item(23, 5).Value = "  +0.49%  "
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = "0.9993"
$ws.Cells.Item(24, 5).Value = "  -0.05%  "
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = "0.1610"
$ws.Cells.Item(25, 5).Value = "  +2.53%  "
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = "162.74"
$ws.Cells.Item(26, 5).Value = "  -0.78%  "
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = "8.914"
$ws.Cells.Item(27, 5).Value = "  -1.22%  "
$ws.Cells.Item(28, 5).Value = "  +1.44%  "
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = "1.361"
$ws.Cells.Item(29, 5).Value = "  -0.36%  "
$ws.Cells.Item(30, 5).Value = "  +1.07%  "
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = "4.440"
$ws.Cells.Item(31, 5).Value = "  +1.13%  "
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = "4.254"
$ws.Cells.Item(32, 5).Value = "  +2.72%  "
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = "0.05182"
$ws.Cells.Item(33, 5).Value = "  -1.87%  "
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = "0.8096"
$ws.Cells.Item(34, 5).Value = "  +12.09%  "
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = "1.932"
$ws.Cells.Item(35, 5).Value = "  -0.17%  "
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = "1.174"
$ws.Cells.Item(36, 5).Value = "  -2.25%  "
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = "2.682"
$ws.Cells.Item(37, 5).Value = "  +0.17%  "
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = "0.01859"
$ws.Cells.Item(38, 5).Value = "  -0.35%  "
$ws.Cells.Item(39, 5).Value = "  -1.07%  "
$ws.Cells.Item(40, 4).Value = "1.167.76"
$ws.Cells.Item(40, 5).Value = "  -5.44%  "
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = "6.218"
$ws.Cells.Item(41, 5).Value = "  +2.25%  "
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = "0.9029"
$ws.Cells.Item(42, 5).Value = "  -0.89%  "
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = "72.82"
$ws.Cells.Item(43, 5).Value = "  -0.70%  "
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = "0.9992"
$ws.Cells.Item(44, 5).Value = "  -0.06%  "
$ws.Cells.Item(45, 2).Value = "Quant"
$ws.Cells.Item(45, 3).Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = "101.98"
$ws.Cells.Item(45, 5).Value = "  -1.63%  "
$ws.Cells.Item(46, 2).Value = "RocketPoolETH"
$ws.Cells.Item(46, 3).Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Cells.Item(46, 4).Value = "2.046.72"
$ws.Cells.Item(46, 5).Value = "  +1.84%  "
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = "0.5173"
$ws.Cells.Item(47, 5).Value = "  -3.07%  "
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = "1.783"
$ws.Cells.Item(48, 5).Value = "  +1.42%  "
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = "9.372"
$ws.Cells.Item(49, 5).Value = "  +1.20%  "
$ws.Cells.Item(50, 5).Value = "  -0.91%  "
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = "7.084"
$ws.Cells.Item(51, 5).Value = "  +0.07%  "
